$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5, shifting rows 5-36 down to 6-37.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the weekly price-record data.
$ws.Cells.Item(5, 1).Value = 7
$ws.Cells.Item(5, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value = "Ñuble"
$ws.Cells.Item(5, 4).Value = 44537
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 100112026
$ws.Cells.Item(5, 7).Value = "Haba"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 60
$ws.Cells.Item(5, 11).Value = 6500
$ws.Cells.Item(5, 12).Value = 7000
$ws.Cells.Item(5, 13).Value = 6750
$ws.Cells.Item(5, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(5, 16).Value = 270
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
